{"js": "// Add a comment \"Nicht informierte Nutzer einf\u00fchren\" anchored on the\n// document title text \"Beschreibung ERM SQL\", and remove the stray\n// \"_GoBack\" bookmark that used to sit at the very end of the body\n// (Word relocates/retires that bookmark once a new edit - the comment -\n// becomes the most recent edit position).\n\nconst body = context.document.body;\n\n// Locate the exact title text so the comment range wraps precisely\n// around it (mirrors selecting the text in the UI before adding a\n// comment, which keeps the comment-reference run after the range end).\nconst titleResults = body.search(\"Beschreibung ERM SQL\", { matchCase: true });\ntitleResults.load(\"items\");\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nlet titleRange;\nif (titleResults.items.length > 0) {\n  titleRange = titleResults.items[0];\n} else {\n  // Fallback: anchor on the title paragraph itself.\n  titleRange = paragraphs.items[0].getRange();\n}\ntitleRange.insertComment(\"Nicht informierte Nutzer einf\u00fchren\");\n\n// Drop the obsolete \"_GoBack\" bookmark left over from the last save.\ncontext.document.deleteBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# Add a comment \"Nicht informierte Nutzer einf\u00fchren\" anchored on the\n# document title text \"Beschreibung ERM SQL\", and remove the stray\n# \"_GoBack\" bookmark that used to sit at the very end of the body\n# (Word relocates/retires that bookmark once a new edit - the comment -\n# becomes the most recent edit position).\n\n$word.UserName = \"Mueller, Kai\"\n\n$d = $word.ActiveDocument\n\n# Find the exact title text so the comment range wraps precisely around\n# it (mirrors selecting the text in the UI before adding a comment,\n# which keeps the comment-reference run after the range end).\n$titleRange = $d.Content\n$found = $titleRange.Find.Execute(\"Beschreibung ERM SQL\")\nif (-not $found) {\n    # Fallback: anchor on the title paragraph itself.\n    $titleRange = $d.Paragraphs(1).Range\n}\n\n$d.Comments.Add($titleRange, \"Nicht informierte Nutzer einf\u00fchren\") | Out-Null\n\n# Drop the obsolete \"_GoBack\" bookmark left over from the last save.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $goBack = $d.Bookmarks.Item(\"_GoBack\")\n    $goBack.Delete()\n}\n"}
